# Apply the "annual_deaths" workbook edits:
#  - Correct the canton name "Geneva" -> "Genève"
#  - Rename the exposure categories "O" -> "Over 75" and "U" -> "Under 75"
#  - Update the sheet view selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix canton name in both halves of the table (Over 75 block row 10, Under 75 block row 37)
$ws.Range("A10").Value = "Genève"
$ws.Range("A37").Value = "Genève"

# Rename category labels
$ws.Range("E2:E28").Value = "Over 75"
$ws.Range("E29:E55").Value = "Under 75"

# Update the active selection shown when the file is opened
$ws.Range("A10").Select()
